$d = $word.ActiveDocument

# --- Q1: merge the fragmented runs (old spell/grammar-check artifacts) into
#         a single run with the same visible text. ---
$p1 = $d.Paragraphs(3).Range
$p1.Find.Execute(
    "Q1. Write a function add_numbers(a, b) that takes two arguments and returns their sum.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q1. Write a function add_numbers(a, b) that takes two arguments and returns their sum.",
    2)

# --- Q2: same kind of run merge. ---
$p2 = $d.Paragraphs(4).Range
$p2.Find.Execute(
    "Q2. Write a function is_even(n) that takes an integer and returns True if the number is even, otherwise False.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q2. Write a function is_even(n) that takes an integer and returns True if the number is even, otherwise False.",
    2)

# --- Q4: same kind of run merge. ---
$p4 = $d.Paragraphs(6).Range
$p4.Find.Execute(
    "Q4. Write a function max_of_two(a, b) that takes two numbers and returns the larger of the two.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q4. Write a function max_of_two(a, b) that takes two numbers and returns the larger of the two.",
    2)

# --- Q5: replace the count_vowels question with the celsius_to_fahrenheit
#         question (content that used to live in Q6), keeping the leading
#         "Q" / "5" split into their own runs. ---
$p5 = $d.Paragraphs(7).Range
$p5.Find.Execute(
    "Q5. Write a function count_vowels(s) that takes a string and returns the number of vowels in the string.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q5. Write a function celsius_to_fahrenheit(c) that takes a temperature in Celsius and returns the temperature in Fahrenheit. (Formula: F = C * 9/5 + 32)",
    2)
$p5start = $d.Paragraphs(7).Range.Start
$p5digit = $d.Range($p5start + 2, $p5start + 3)
$p5digit.Font.Bold = 1
$p5digit.Font.Bold = 0

# --- Q6: replace the celsius_to_fahrenheit question with the
#         multiplication_table question (content that used to live in Q7),
#         again keeping "Q" / "6" as their own runs. ---
$p6 = $d.Paragraphs(8).Range
$p6.Find.Execute(
    "Q6. Write a function celsius_to_fahrenheit(c) that takes a temperature in Celsius and returns the temperature in Fahrenheit. (Formula: F = C * 9/5 + 32)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q6. Write a function multiplication_table(n) that takes an integer and prints its multiplication table from 1 to 10.",
    2)
$p6start = $d.Paragraphs(8).Range.Start
$p6digit = $d.Range($p6start + 2, $p6start + 3)
$p6digit.Font.Bold = 1
$p6digit.Font.Bold = 0

# --- Q7: the whole paragraph (its former multiplication_table question) is
#         removed outright now that its text lives in Q6. ---
$p7 = $d.Paragraphs(9).Range
$p7.Delete()
